$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.378.82'
$ws.Range('E2').Value = '  -1.17%  '

# Row 3
$ws.Range('D3').Value = '1.816.90'
$ws.Range('E3').Value = '  -3.12%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.99%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.99'
$ws.Range('E5').Value = '  -1.32%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.85%  '

# Row 7
$ws.Range('E7').Value = '  -2.06%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3798'
$ws.Range('E8').Value = '  -3.77%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.63'
$ws.Range('E9').Value = '  -0.80%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07845'
$ws.Range('E10').Value = '  -1.97%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9555'
$ws.Range('E11').Value = '  -5.12%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.94'
$ws.Range('E12').Value = '  -4.24%  '

# Row 13
$ws.Range('D13').Value = '1.833.61'
$ws.Range('E13').Value = '  -2.45%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.825'
$ws.Range('E14').Value = '  -2.63%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.054'
$ws.Range('E15').Value = '  -2.69%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.005'
$ws.Range('E16').Value = '  -0.93%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.04'
$ws.Range('E17').Value = '  +0.18%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06573'
$ws.Range('E18').Value = '  -2.38%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001018'
$ws.Range('E19').Value = '  -2.68%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.11'
$ws.Range('E20').Value = '  -1.11%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  -0.81%  '

# Row 22
$ws.Range('D22').Value = '27.365.47'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.267'
$ws.Range('E23').Value = '  -3.92%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.79'
$ws.Range('E24').Value = '  -1.59%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.263'
$ws.Range('E25').Value = '  -1.77%  '

# Row 26
$ws.Range('D26').Value = '2.047.82'
$ws.Range('E26').Value = '  -2.55%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.89'
$ws.Range('E27').Value = '  -2.56%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.27'
$ws.Range('E28').Value = '  -2.44%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.039'
$ws.Range('E29').Value = '  -5.27%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.252'
$ws.Range('E30').Value = '  -4.05%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.31'
$ws.Range('E31').Value = '  -3.58%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09273'
$ws.Range('E32').Value = '  -2.01%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9298'
$ws.Range('E33').Value = '  -5.14%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.564'
$ws.Range('E34').Value = '  -1.72%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.201'
$ws.Range('E35').Value = '  -2.46%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.312'
$ws.Range('E36').Value = '  -2.61%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05895'
$ws.Range('E37').Value = '  -2.73%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02174'
$ws.Range('E38').Value = '  -3.04%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.084'
$ws.Range('E39').Value = '  -3.34%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.003'
$ws.Range('E40').Value = '  -0.81%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.138'
$ws.Range('E41').Value = '  -5.35%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5733'
$ws.Range('E42').Value = '  -4.00%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1809'
$ws.Range('E43').Value = '  -3.78%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.898'
$ws.Range('E44').Value = '  -4.68%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.275'
$ws.Range('E45').Value = '  +1.96%  '

# Row 46
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5382'
$ws.Range('E46').Value = '  -4.72%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.78'
$ws.Range('E47').Value = '  -3.37%  '

# Row 48
$ws.Range('E48').Value = '  -3.65%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06568'
$ws.Range('E49').Value = '  -2.90%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '109.47'
$ws.Range('E50').Value = '  -2.43%  '

# Row 51
$ws.Range('E51').Value = '  -33.66%  '
